$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("2024-09-02", 1092.650024414062, 715.0499877929688, 251.3500061035156, 490.5, 1505.25, 29322.30020141602, 0, 248.8582410022039),
  @("2024-09-03", 1085.099975585938, 710.7999877929688, 251, 488.8500061035156, 1509, 29201.39971923828, -0.004123158188384413, 247.8321591080687),
  @("2024-09-04", 1080.449951171875, 722.4000244140625, 250.5, 484.1499938964844, 1488.099975585938, 29200.89978027344, -0.00001712037675078922, 247.8279161281338),
  @("2024-09-05", 1069.150024414062, 733.8499755859375, 251.1499938964844, 495.6499938964844, 1447.599975585938, 29262.3996887207, 0.002106096350113556, 248.3498655977476),
  @("2024-09-06", 1049.349975585938, 718.9000244140625, 247.8000030517578, 483, 1418.050048828125, 28702.20024108887, -0.01914400232349252, 243.5954551937052),
  @("2024-09-09", 1038.699951171875, 700.1500244140625, 243.8999938964844, 474.75, 1411.849975585938, 28242.6496887207, -0.01601098691069303, 239.6952515490945),
  @("2024-09-10", 1035.800048828125, 713.4000244140625, 248.25, 478.7999877929688, 1424.449951171875, 28522.85046386719, 0.009921192885042528, 242.0733143733418),
  @("2024-09-11", 976.2999877929688, 725.4000244140625, 241.5500030517578, 472.2000122070312, 1399.599975585938, 27922.25028991699, -0.02105680758348599, 236.9760231714857),
  @("2024-09-12", 986.1500244140625, 726.0499877929688, 246.1499938964844, 479.8500061035156, 1403.150024414062, 28182.30001831055, 0.009313351384414074, 239.1830641449628),
  @("2024-09-13", 992.0999755859375, 724.25, 245.6499938964844, 485.3999938964844, 1410.949951171875, 28258.59951782227, 0.002707355306775728, 239.8306176829666),
  @("2024-09-16", 988.4000244140625, 733.6500244140625, 243.8000030517578, 489.9500122070312, 1404.550048828125, 28313.45072937012, 0.00194104493795801, 240.2961396893874),
  @("2024-09-17", 974.9500122070312, 745.4000244140625, 240.8000030517578, 482.2999877929688, 1400.25, 28196.30033874512, -0.004137623200533361, 239.30188480681),
  @("2024-09-18", 962.0499877929688, 717.5499877929688, 235.9499969482422, 471.75, 1391.300048828125, 27572.89979553223, -0.02210930284198537, 234.0110869649583),
  @("2024-09-19", 967, 728.5, 237.5500030517578, 459.9500122070312, 1374.150024414062, 27641.0502166748, 0.002471645044516532, 234.5894793084172),
  @("2024-09-20", 970.8499755859375, 748.3499755859375, 237.8500061035156, 466.2999877929688, 1380.550048828125, 27960.69967651367, 0.01156430227264067, 237.3023429571211),
  @("2024-09-23", 971.7999877929688, 750.2000122070312, 236.4499969482422, 471.1499938964844, 1375.400024414062, 27984.94996643066, 0.0008672991090191444, 237.508155067736),
  @("2024-09-24", 977.2999877929688, 735.9000244140625, 237.3000030517578, 476.7000122070312, 1363.699951171875, 27912.20024108887, -0.002599601765558408, 236.8907284484874),
  @("2024-09-25", 963.5999755859375, 730.0499877929688, 238.3500061035156, 473.7000122070312, 1365.400024414062, 27741.79995727539, -0.006104867489544392, 235.4445419418077),
  @("2024-09-26", 993.1500244140625, 744.0999755859375, 241.1999969482422, 471.75, 1329.949951171875, 28111.79975891113, 0.01333726730801793, 238.5847287338994),
  @("2024-09-27", 993, 735.4500122070312, 239.5500030517578, 497.2999877929688, 1392.199951171875, 28306.40000915527, 0.006922368966521059, 240.2363002559728)
)

$startRow = 654
$rowIdx = $startRow
foreach ($r in $data) {
    $ws.Cells.Item($rowIdx, 1).NumberFormat = "@"
    $ws.Cells.Item($rowIdx, 1).Value = $r[0]
    $ws.Cells.Item($rowIdx, 1).ClearFormats()

    $ws.Cells.Item($rowIdx, 3).Value = $r[1]
    $ws.Cells.Item($rowIdx, 4).Value = $r[2]
    $ws.Cells.Item($rowIdx, 5).Value = $r[3]
    $ws.Cells.Item($rowIdx, 6).Value = $r[4]
    $ws.Cells.Item($rowIdx, 7).Value = $r[5]
    $ws.Cells.Item($rowIdx, 8).Value = $r[6]
    $ws.Cells.Item($rowIdx, 9).Value = $r[7]
    $ws.Cells.Item($rowIdx, 10).Value = $r[8]

    $rowIdx = $rowIdx + 1
}
